$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 22, shifting the existing data
# (rows 22-44) down to rows 24-46.
$ws.Rows.Item(22).Resize(2).Insert()

# Populate the two new rows (22 and 23) with the new weekly records.
$newRows = @(
    @{ Row = 22; D = 44546; L = "Primera"; M = 100;  N = 3800; O = 3800; P = 3800; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares"; S = 1900; T = 2 },
    @{ Row = 23; D = 44546; L = "Segunda"; M = 150;  N = 3000; O = 3000; P = 3000; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares"; S = 1500; T = 2 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100101001
    $ws.Cells.Item($row, 10).Value = "Arándano (blue)"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
